$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 19905.188
$ws.Range("J32").Value = 20458.666
$ws.Range("L32").Value = 20458.666
$ws.Range("N32").Value = -21110.666
$ws.Range("H33").Value = 675.9048
$ws.Range("I33").Value = 615.75
$ws.Range("K33").Value = 615.75
$ws.Range("M33").Value = -386.75
$ws.Range("H51").Value = 20882.572
$ws.Range("J51").Value = 6031.3335
$ws.Range("L51").Value = 6031.3335
$ws.Range("N51").Value = -6999.3335
$ws.Range("H53").Value = 509.58334
$ws.Range("I53").Value = 267.875
$ws.Range("K53").Value = 267.875
$ws.Range("M53").Value = 369.125
$ws.Range("H132").Value = 10788.333
$ws.Range("I132").Value = 9011.875
$ws.Range("K132").Value = 27035.625
$ws.Range("M132").Value = -24505.625
$ws.Range("H137").Value = 8374.076999999999
$ws.Range("I137").Value = 5514.875
$ws.Range("K137").Value = 16544.625
$ws.Range("M137").Value = -13994.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 8244.846
$ws.Range("I102").Value = 8244.846
$ws.Range("K102").Value = 8244.846
$ws.Range("M102").Value = -6622.846
$ws.Range("H122").Value = 3775.6487
$ws.Range("I122").Value = 3368.3572
$ws.Range("K122").Value = 10105.0716
$ws.Range("M122").Value = -7655.071599999999
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 45000
$ws.Range("J130").Value = 45000
$ws.Range("L130").Value = 45000
$ws.Range("N130").Value = -55040
$ws.Range("H131").Value = 63999.5
$ws.Range("J131").Value = 63999.5
$ws.Range("L131").Value = 63999.5
$ws.Range("N131").Value = -74079.5
$ws.Range("H132").Value = 44945.59
$ws.Range("I132").Value = 1868.4286
$ws.Range("K132").Value = 5605.2858
$ws.Range("M132").Value = -3075.2858
$ws.Range("H140").Value = 41999
$ws.Range("J140").Value = 41999
$ws.Range("L140").Value = 41999
$ws.Range("N140").Value = -52359

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2715.2
$ws.Range("I22").Value = 1289.6
$ws.Range("J22").Value = 5566.4
$ws.Range("K22").Value = 1289.6
$ws.Range("L22").Value = 5566.4
$ws.Range("M22").Value = -1116.6
$ws.Range("N22").Value = -5912.4
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H119").Value = 87708.2
$ws.Range("J119").Value = 87708.2
$ws.Range("L119").Value = 87708.2
$ws.Range("N119").Value = -97384.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5047.5
$ws.Range("J16").Value = 6724.5
$ws.Range("L16").Value = 6724.5
$ws.Range("N16").Value = -7298.5
$ws.Range("H31").Value = 7068.4614
$ws.Range("I31").Value = 2456.5715
$ws.Range("J31").Value = 12449
$ws.Range("K31").Value = 2456.5715
$ws.Range("L31").Value = 12449
$ws.Range("M31").Value = -2161.5715
$ws.Range("N31").Value = -13039
$ws.Range("H34").Value = 7068.4614
$ws.Range("I34").Value = 2456.5715
$ws.Range("J34").Value = 12449
$ws.Range("K34").Value = 2456.5715
$ws.Range("L34").Value = 12449
$ws.Range("M34").Value = -2254.5715
$ws.Range("N34").Value = -12853
$ws.Range("H43").Value = 16252
$ws.Range("J43").Value = 16252
$ws.Range("L43").Value = 16252
$ws.Range("N43").Value = -16620
$ws.Range("H58").Value = 12339.647
$ws.Range("I58").Value = 9617.23
$ws.Range("K58").Value = 9617.23
$ws.Range("M58").Value = -9414.23
$ws.Range("H95").Value = 16849.834
$ws.Range("J95").Value = 16849.834
$ws.Range("L95").Value = 16849.834
$ws.Range("N95").Value = -22341.834
$ws.Range("H99").Value = 4818
$ws.Range("I99").Value = 4715.2
$ws.Range("K99").Value = 4715.2
$ws.Range("M99").Value = -3217.2
$ws.Range("H101").Value = 16252
$ws.Range("J101").Value = 16252
$ws.Range("L101").Value = 16252
$ws.Range("N101").Value = -22742
$ws.Range("H113").Value = 5047.5
$ws.Range("J113").Value = 6724.5
$ws.Range("L113").Value = 6724.5
$ws.Range("N113").Value = -11064.5
$ws.Range("H126").Value = 4818
$ws.Range("I126").Value = 4715.2
$ws.Range("K126").Value = 14145.6
$ws.Range("M126").Value = -11675.6
$ws.Range("H132").Value = 7081.129
$ws.Range("I132").Value = 6367.222
$ws.Range("K132").Value = 19101.666
$ws.Range("M132").Value = -16571.666
$ws.Range("H136").Value = 12339.647
$ws.Range("I136").Value = 9617.23
$ws.Range("K136").Value = 28851.69
$ws.Range("M136").Value = -26301.69

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 320717.88
$ws.Range("I68").Value = 7499.3335
$ws.Range("J68").Value = 392999.06
$ws.Range("K68").Value = 22498.0005
$ws.Range("L68").Value = 1178997.18
$ws.Range("M68").Value = -21687.0005
$ws.Range("N68").Value = -1180619.18
$ws.Range("H71").Value = 320717.88
$ws.Range("I71").Value = 7499.3335
$ws.Range("J71").Value = 392999.06
$ws.Range("K71").Value = 67494.0015
$ws.Range("L71").Value = 3536991.54
$ws.Range("M71").Value = -63438.0015
$ws.Range("N71").Value = -3545103.54
$ws.Range("H75").Value = 1030.4286
$ws.Range("I75").Value = 1153.25
$ws.Range("K75").Value = 3459.75
$ws.Range("M75").Value = -2461.75
$ws.Range("H78").Value = 1030.4286
$ws.Range("I78").Value = 1153.25
$ws.Range("K78").Value = 10379.25
$ws.Range("M78").Value = -5387.25
$ws.Range("H107").Value = 5022.2666
$ws.Range("J107").Value = 5542.926
$ws.Range("L107").Value = 16628.778
$ws.Range("N107").Value = -20468.778
$ws.Range("H113").Value = 500749.8
$ws.Range("I113").Value = 340.5
$ws.Range("J113").Value = 625852.1
$ws.Range("K113").Value = 1021.5
$ws.Range("L113").Value = 1877556.3
$ws.Range("M113").Value = 1148.5
$ws.Range("N113").Value = -1881896.3
$ws.Range("H114").Value = 1269
$ws.Range("J114").Value = 1300
$ws.Range("L114").Value = 3900
$ws.Range("N114").Value = -10408
$ws.Range("H141").Value = 7030
$ws.Range("I141").Value = 7030
$ws.Range("K141").Value = 21090
$ws.Range("M141").Value = -15910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6687
$ws.Range("J102").Value = 7607.8887
$ws.Range("L102").Value = 7607.8887
$ws.Range("N102").Value = -10851.8887
$ws.Range("H123").Value = 43499.5
$ws.Range("J123").Value = 43499.5
$ws.Range("L123").Value = 43499.5
$ws.Range("N123").Value = -48399.5
$ws.Range("H127").Value = 75749.25
$ws.Range("J127").Value = 75749.25
$ws.Range("L127").Value = 75749.25
$ws.Range("N127").Value = -85669.25
$ws.Range("H129").Value = 29888.666
$ws.Range("J129").Value = 29888.666
$ws.Range("L129").Value = 29888.666
$ws.Range("N129").Value = -39888.666
$ws.Range("H131").Value = 34500
$ws.Range("J131").Value = 34500
$ws.Range("L131").Value = 34500
$ws.Range("N131").Value = -44580
$ws.Range("H132").Value = 1296.9
$ws.Range("I132").Value = 1083.625
$ws.Range("K132").Value = 3250.875
$ws.Range("M132").Value = -720.875
$ws.Range("H134").Value = 47725
$ws.Range("J134").Value = 47725
$ws.Range("L134").Value = 143175
$ws.Range("N134").Value = -148245

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1132.6666
$ws.Range("I22").Value = 1099.5
$ws.Range("J22").Value = 1199
$ws.Range("K22").Value = 1099.5
$ws.Range("L22").Value = 1199
$ws.Range("M22").Value = -804.5
$ws.Range("N22").Value = -1789
$ws.Range("H27").Value = 1132.6666
$ws.Range("I27").Value = 1099.5
$ws.Range("J27").Value = 1199
$ws.Range("K27").Value = 1099.5
$ws.Range("L27").Value = 1199
$ws.Range("M27").Value = -992.5
$ws.Range("N27").Value = -1413
$ws.Range("H68").Value = 1918.7059
$ws.Range("I68").Value = 1838.625
$ws.Range("J68").Value = 3200
$ws.Range("K68").Value = 1838.625
$ws.Range("L68").Value = 3200
$ws.Range("M68").Value = -1089.625
$ws.Range("N68").Value = -4698
$ws.Range("H71").Value = 1918.7059
$ws.Range("I71").Value = 1838.625
$ws.Range("J71").Value = 3200
$ws.Range("K71").Value = 9193.125
$ws.Range("L71").Value = 16000
$ws.Range("M71").Value = -5449.125
$ws.Range("N71").Value = -23488
$ws.Range("H103").Value = 30049.25
$ws.Range("J103").Value = 30049.25
$ws.Range("L103").Value = 30049.25
$ws.Range("N103").Value = -32393.25
$ws.Range("H135").Value = 22714.5
$ws.Range("J135").Value = 22714.5
$ws.Range("L135").Value = 22714.5
$ws.Range("N135").Value = -32854.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7400.9287
$ws.Range("I132").Value = 4877
$ws.Range("K132").Value = 14631
$ws.Range("M132").Value = -14631
$ws.Range("H136").Value = 68462.89
$ws.Range("J136").Value = 28577
$ws.Range("L136").Value = 85731
$ws.Range("N136").Value = -90831
